$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.546.78"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.045.30"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'231.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -12.00%  "
$ws.Range("D6").Value = "'0.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'55.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "'0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'56.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "'0.0748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "2.342.43"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").Value = "'14.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'20.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.05%  "
$ws.Range("D16").Value = "'0.761"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'5.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "2.035.20"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "36.773.25"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "'67.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "'5.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.75%  "
$ws.Range("D22").Value = "0.0₃0796"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'220.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.20%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "'2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.95%  "
$ws.Range("D27").Value = "'162.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'8.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'0.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'18.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").Value = "'1.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "'4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "'2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("D35").Value = "'0.0603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'4.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("D40").Value = "'3.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").Value = "'2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("D42").Value = "1.475.72"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +42.97%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0932"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'93.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.90%  "
$ws.Range("D46").Value = "'0.0204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "'15.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "'1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'2.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'6.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.75%  "
